# Apply cryptos list update (prices + volume%) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.230.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.293.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.95%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("E11").Value = "  -1.15%  "

$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.965"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.641.11"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.296.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.247.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "

$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +29.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("E25").Value = "  -4.41%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "

$ws.Range("E28").Value = "  +2.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0875"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.47%  "

$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0356"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.60%  "

$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.23%  "

$ws.Range("E41").Value = "  +3.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "94.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.97%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "116.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("E47").Value = "  -3.88%  "

$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.605.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.11%  "
